# daily auto push: 2026-02-15 09:43 UTC
# Insert one new data row at row 799 (pushing the existing rows 799-840
# down to 800-841) and populate the new row with the day's first
# observation: 2026/02/15, 日, time=16, ranking=29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 799..840 down by one row, opening up a blank row 799.
$ws.Rows.Item(799).Insert()

# Column A holds plain text dates (e.g. "2026/12/29"), not real dates.
# A bare assignment of a date-shaped string gets auto-converted to a
# date serial by Excel's smart typing, so force text entry with a
# leading apostrophe and then strip the resulting quote-prefix style so
# the cell matches the plain (unstyled) text cells used everywhere else
# in this column.
$ws.Cells.Item(799, 1).Value = "'2026/02/15"
$ws.Cells.Item(799, 1).ClearFormats()

$ws.Cells.Item(799, 2).Value = "日"
$ws.Cells.Item(799, 3).Value = 16
$ws.Cells.Item(799, 4).Value = 29
